$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures scraped on 2024-02-04 (plus the MultiversX /
# RocketPoolETH row swap at the bottom of the table).
#
# Column D ("Price") cells are stored as plain text in the workbook (e.g. "299.41",
# "2.295.33"). Assigning a numeric-looking string straight to .Value would make Excel
# reinterpret it as a number (and reformat/round it), so for every Price cell we force
# the cell to Text format first, write the literal string, then restore the original
# (unstyled) cell style so no stray formatting is introduced.
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "42.934.62"
$ws.Range("E2").Value = "  -0.59%  "
Set-TextValue "D3" "2.295.33"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue "D5" "299.41"
$ws.Range("E5").Value = "  -1.16%  "
Set-TextValue "D6" "97.35"
$ws.Range("E6").Value = "  -2.16%  "
Set-TextValue "D7" "0.515"
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("E8").Value = "  -0.04%  "
Set-TextValue "D9" "0.507"
$ws.Range("E9").Value = "  -1.93%  "
Set-TextValue "D10" "35.73"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("E12").Value = "  +0.69%  "
Set-TextValue "D13" "17.66"
$ws.Range("E13").Value = "  +0.04%  "
Set-TextValue "D14" "6.77"
$ws.Range("E14").Value = "  -2.18%  "
Set-TextValue "D15" "2.654.81"
$ws.Range("E15").Value = "  -1.11%  "
Set-TextValue "D16" "2.315.71"
$ws.Range("E16").Value = "  -0.07%  "
Set-TextValue "D17" "0.776"
$ws.Range("E17").Value = "  -2.59%  "
Set-TextValue "D18" "42.877.92"
$ws.Range("E18").Value = "  -0.50%  "
Set-TextValue "D19" "12.59"
$ws.Range("E19").Value = "  -4.75%  "
Set-TextValue "D20" "0.0₃0907"
$ws.Range("E20").Value = "  -0.57%  "
Set-TextValue "D21" "6.09"
$ws.Range("E21").Value = "  -2.42%  "
Set-TextValue "D22" "68.05"
$ws.Range("E22").Value = "  -0.12%  "
Set-TextValue "D23" "241.66"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("E27").Value = "  -0.39%  "
Set-TextValue "D28" "25.08"
$ws.Range("E28").Value = "  -1.68%  "
Set-TextValue "D29" "166.51"
$ws.Range("E29").Value = "  -0.89%  "
Set-TextValue "D31" "9.05"
$ws.Range("E31").Value = "  -1.63%  "
Set-TextValue "D32" "32.83"
$ws.Range("E32").Value = "  -4.17%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  -3.43%  "
Set-TextValue "D35" "4.69"
$ws.Range("E35").Value = "  -0.97%  "
Set-TextValue "D36" "17.61"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("E39").Value = "  -2.73%  "
$ws.Range("E40").Value = "  -3.53%  "
Set-TextValue "D41" "2.75"
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("E42").Value = "  +0.16%  "
Set-TextValue "D43" "2.000.85"
$ws.Range("E43").Value = "  +0.38%  "
Set-TextValue "D44" "0.0286"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("E45").Value = "  -3.89%  "
Set-TextValue "D46" "10.14"
$ws.Range("E46").Value = "  +0.38%  "
Set-TextValue "D47" "17.37"
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("E48").Value = "  -3.83%  "
Set-TextValue "D49" "2.92"
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D50" "2.524.24"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D51" "53.32"
$ws.Range("E51").Value = "  -3.13%  "
